$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-01 12:50:08"

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
